# Fruta / hortaliza, semanal
# A new weekly price observation (2021-11-09) is inserted into the
# Betarraga / Feria Lagunitas de Puerto Montt dataset, right above the
# most recent existing row (which held the 2021-10-22 observation).
# This pushes the existing rows 199-206 down to 200-207.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 199; Excel shifts rows 199-206 down to 200-207,
# carrying over formatting (incl. the date style on column D) from the
# row above.
$ws.Rows.Item(199).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Cells.Item(199, 1).Value  = 4
$ws.Cells.Item(199, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(199, 3).Value  = "Los Lagos"
$ws.Cells.Item(199, 4).Value  = 44509
$ws.Cells.Item(199, 5).Value  = 10
$ws.Cells.Item(199, 6).Value  = 100114014
$ws.Cells.Item(199, 7).Value  = "Betarraga"
$ws.Cells.Item(199, 8).Value  = "Sin especificar"
$ws.Cells.Item(199, 9).Value  = "Primera"
$ws.Cells.Item(199, 10).Value = 1200
$ws.Cells.Item(199, 11).Value = 1000
$ws.Cells.Item(199, 12).Value = 1200
$ws.Cells.Item(199, 13).Value = 1100
$ws.Cells.Item(199, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(199, 15).Value = "Región del Maule"
$ws.Cells.Item(199, 16).Value = 220
$ws.Cells.Item(199, 17).Value = 5
$ws.Cells.Item(199, 18).Value = "Hortaliza"
